# Trade #20 closed at 2026-02-17 23:56:58 - unknown UNKNOWN +0.000%
#
# Updates the "Summary", "Strategy Status", "All Trades" and
# "MarketMaking" sheets of the live trading results workbook to reflect
# Trade #20 (row 21 on the trade-log sheets) transitioning from OPEN to
# CLOSED.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1500.76   # Current Capital
$summary.Range("B4").Value = 0.76      # Total P&L $
$summary.Range("B5").Value = 0.76      # Total P&L %
$summary.Range("B6").Value = 20        # Total Trades
$summary.Range("B7").Value = 11        # Winning Trades
$summary.Range("B9").Value = 55        # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet - MarketMaking row (row 6)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C6").Value = 100.76     # Capital
$status.Range("D6").Value = 20         # Trades
$status.Range("E6").Value = 0.76       # P&L $
$status.Range("F6").Value = 0.76       # P&L %
$status.Range("G6").Value = 55         # Win Rate %

# ---------------------------------------------------------------------
# All Trades + MarketMaking trade-log sheets - Trade #20 (row 21)
# ---------------------------------------------------------------------
$tradeSheets = @("All Trades", "MarketMaking")
foreach ($sheetName in $tradeSheets) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("G21").Value = 0.45          # Exit Price
    $ws.Range("H21").Value = "CLOSED"      # Status
    $ws.Range("I21").Value = 4400          # P&L %
    $ws.Range("J21").Value = 0.44          # P&L $
    $ws.Range("K21").Value = 100.76        # Capital After
    $ws.Range("P21").Value = "early_exit"  # Exit Reason
    $ws.Range("Q21").Value = 2.33          # Duration (min)
}
